$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bill_items")

# Make room for a new "item_description" column at column C by shifting the
# existing columns C:G one position to the right (work right-to-left so the
# source of each copy hasn't been overwritten yet). Copy (rather than plain
# value assignment) so that cell formatting/styles move along with the data.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("D1").Copy($ws.Range("E1"))
$ws.Range("C1").Copy($ws.Range("D1"))

$ws.Range("F2").Copy($ws.Range("G2"))
$ws.Range("E2").Copy($ws.Range("F2"))
$ws.Range("D2").Copy($ws.Range("E2"))
$ws.Range("C2").Copy($ws.Range("D2"))

# Fill in the new "item_description" column with its header and sample value
$ws.Range("C1").Value2 = "item_description"
$ws.Range("C1").Style = "Normal"
$ws.Range("C2").Value2 = "test bill item"
$ws.Range("C2").Style = "Normal"

# Widen the new column to fit its header/content
$ws.Columns("C:C").ColumnWidth = 15.17

# Put the selection/active cell on the newly filled-in data cell, matching
# where the user last made their edit
$ws.Range("C2").Select() | Out-Null
